# "Adicao das variaveis inteiras e transformacao do problema em MINLP"
#
# The generation-options table gains:
#   - two new integer-variable columns ("existentes" / "novas") inserted
#     right after "Tipo" and before "Capacidade";
#   - three new trailing columns ("consumo" / "limite_comb" / "Emissao")
#     used by the new MINLP fuel-consumption / emission constraints.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Capacidade..OeM vari" block (old B:G) two columns to
# the right (new D:I) by inserting two blank columns at B:C.
$ws.Range("B:C").Insert()

$ws.Range("B1").Value = "existentes"
$ws.Range("C1").Value = "novas"

# existentes / novas per plant type (Oleo, Carvao, Gas, Biomassa, Eolica, PV)
$existentes = 1, 1, 1, 0, 1, 0
$novas      = 0, 1, 2, 0, 2, 2

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $existentes[$i]
    $ws.Range("C$row").Value = $novas[$i]
}

# New trailing columns for the MINLP fuel/emission limits, all zero-initialised.
$ws.Range("J1").Value = "consumo"
$ws.Range("K1").Value = "limite_comb"
$ws.Range("L1").Value = "Emissao"

for ($row = 2; $row -le 7; $row++) {
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 0
    $ws.Range("L$row").Value = 0
}

$ws.Range("C2").Select()
